# Staging.Programme.xlsx - field list reorder.
#
# The original sheet1 header row (row 2) mapped shared strings as:
#   A2=Programme_ID  B2=Code  C2=LongName  D2=BusinessKey
#   E2=ShortName     F2=TextDescription    G2=ProgrammeSiteName
#
# After the edit the fields were reshuffled (BusinessKey promoted to the
# front, ShortName/TextDescription/ProgrammeSiteName rotated) so the
# row now reads:
#   A2=BusinessKey   B2=Code  C2=LongName  D2=Programme_ID
#   E2=ProgrammeSiteName      F2=ShortName G2=TextDescription
# B2 ("Code") and C2 ("LongName") are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "BusinessKey"
$ws.Range("D2").Value = "Programme_ID"
$ws.Range("E2").Value = "ProgrammeSiteName"
$ws.Range("F2").Value = "ShortName"
$ws.Range("G2").Value = "TextDescription"

# Best-effort: restore the saved window size recorded in the workbook view.
# (Read-only/unsupported through plain COM automation on some hosts, so this
# is wrapped defensively and is a no-op if unsupported.)
try { $excel.ActiveWindow.Width = 28800 } catch {}
try { $excel.ActiveWindow.Height = 12585 } catch {}

# Best-effort: the sheet's internal VBA CodeName moved from Sheet41 to
# Sheet43 in the authored file. CodeName is normally read-only from
# external automation (exactly like real Excel), so this is attempted but
# tolerated if the host rejects/ignores it.
try { $ws.CodeName = "Sheet43" } catch {}

# Note: columns B:G also lost their manually recorded bestFit widths in the
# authored file (the <col> customizations disappeared entirely, reverting
# to the sheet's default width). There is no Excel object-model call that
# clears a column's width back to "never set" - ColumnWidth/AutoFit always
# (re)write an explicit, custom width - so attempting that here would trade
# one set of custom widths for a different, equally-wrong set. Leaving the
# original widths in place is the more faithful no-op.
